# Refresh crypto price/volume snapshot (columns D and E) for rows 2-51,
# matching the scheduled GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.897.48'
$ws.Range("E2").Value = '  +0.09%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.889.28'
$ws.Range("E3").Value = '  -0.13%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7653'
$ws.Range("E5").Value = '  -1.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.79'
$ws.Range("E6").Value = '  -0.59%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3130'
$ws.Range("E8").Value = '  -0.36%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.68'
$ws.Range("E9").Value = '  +1.43%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07170'
$ws.Range("E10").Value = '  -3.09%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08519'
$ws.Range("E11").Value = '  +4.53%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7624'
$ws.Range("E12").Value = '  -0.52%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.916.54'
$ws.Range("E13").Value = '  +2.32%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.375'
$ws.Range("E14").Value = '  -1.95%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '93.82'
$ws.Range("E15").Value = '  +1.41%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.134'
$ws.Range("E16").Value = '  -1.30%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.816.81'
$ws.Range("E17").Value = '  -0.15%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.78'
$ws.Range("E18").Value = '  -1.24%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.03'
$ws.Range("E19").Value = '  -0.51%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007802'
$ws.Range("E20").Value = '  -0.88%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.0000'
$ws.Range("E21").Value = '  +0.01%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.043'
$ws.Range("E22").Value = '  -1.00%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.110.17'
$ws.Range("E23").Value = '  +0.09%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9998'
$ws.Range("E24").Value = '  -0.06%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1620'
$ws.Range("E25").Value = '  +3.34%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.402'
$ws.Range("E26").Value = '  -0.35%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.13'
$ws.Range("E27").Value = '  -0.06%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.76'
$ws.Range("E28").Value = '  -0.25%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.036'
$ws.Range("E29").Value = '  -0.24%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.491'
$ws.Range("E30").Value = '  +2.72%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.529'
$ws.Range("E31").Value = '  -1.11%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.490'
$ws.Range("E32").Value = '  -0.23%  '

$ws.Range("E33").Value = '  -0.15%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05429'
$ws.Range("E34").Value = '  -3.17%  '

$ws.Range("E35").Value = '  -0.56%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7437'
$ws.Range("E36").Value = '  -1.92%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9999'
$ws.Range("E37").Value = '  +0.15%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.693'
$ws.Range("E38").Value = '  +1.70%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01948'
$ws.Range("E39").Value = '  +0.54%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.780'
$ws.Range("E40").Value = '  -0.30%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4467'
$ws.Range("E41").Value = '  +0.14%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.101.12'
$ws.Range("E42").Value = '  -4.34%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '73.13'
$ws.Range("E43").Value = '  -1.74%  '

$ws.Range("E44").Value = '  +1.60%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8521'
$ws.Range("E45").Value = '  -0.31%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.0000'
$ws.Range("E46").Value = '  -0.04%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.89'
$ws.Range("E47").Value = '  +0.95%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.870'
$ws.Range("E48").Value = '  -1.92%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.653'
$ws.Range("E49").Value = '  +1.87%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.996'
$ws.Range("E50").Value = '  -4.03%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.014.35'
$ws.Range("E51").Value = '  -0.81%  '
